$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet (placing the copy right
# after it) so the current quarter's figures are preserved unchanged under
# the same sheet name. The original sheet object is then repurposed to hold
# the brand-new "2022-Q3" figures. This reproduces the target layout:
#   总计 (sheetId 1) / 2022-Q3 (sheetId 2) / 2022-Q2 (sheetId 3)
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy([System.Reflection.Missing]::Value, $q2)
$q2copy = $wb.Worksheets.Item(3)

$q2.Name = "2022-Q3"
$q2copy.Name = "2022-Q2"

$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# Step 2: update the "总计" (totals) summary sheet.
# The row that used to hold the 2022-Q2 totals moves down to row 3 (with its
# sequence index bumped from 0 to 1); row 2 is rewritten with the new
# 2022-Q3 totals.
# ---------------------------------------------------------------------------
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 3
$total.Cells.Item(3, 4).Value = 0.02

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 1.14

# Row 3's A-cell must keep the same bordered/centered index-column style as
# row 2 (style index 2 in the original workbook).
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item(3, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 3: replace the "2022-Q3" sheet contents with the new fund-holding
# data. The column-B..G values are fund codes / percentages that must stay
# text (e.g. fund code "001239" should not collapse to the number 1239), so
# they are forced to text via a temporary "@" number format, then the style
# is reset to "Normal" so no stray cell format is left behind.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $q3.Cells.Item(1, 2) "基金代码"
Set-TextValue $q3.Cells.Item(1, 3) "基金名称"
Set-TextValue $q3.Cells.Item(1, 4) "基金规模"
Set-TextValue $q3.Cells.Item(1, 5) "股票总仓位"
Set-TextValue $q3.Cells.Item(1, 6) "仓位占比"
Set-TextValue $q3.Cells.Item(1, 7) "持有市值(亿元)"
Set-TextValue $q3.Cells.Item(1, 8) "仓位排名"

$q3.Cells.Item(2, 1).Value = 0
Set-TextValue $q3.Cells.Item(2, 2) "590003"
Set-TextValue $q3.Cells.Item(2, 3) "中邮核心优势灵活配置混合"
Set-TextValue $q3.Cells.Item(2, 4) "18.95"
Set-TextValue $q3.Cells.Item(2, 5) "79.43"
Set-TextValue $q3.Cells.Item(2, 6) "4.90"
Set-TextValue $q3.Cells.Item(2, 7) "0.9286"
$q3.Cells.Item(2, 8).Value = 8

$q3.Cells.Item(3, 1).Value = 1
Set-TextValue $q3.Cells.Item(3, 2) "001239"
Set-TextValue $q3.Cells.Item(3, 3) "长盛国企改革主题灵活配置混合"
Set-TextValue $q3.Cells.Item(3, 4) "4.46"
Set-TextValue $q3.Cells.Item(3, 5) "90.97"
Set-TextValue $q3.Cells.Item(3, 6) "4.75"
Set-TextValue $q3.Cells.Item(3, 7) "0.2118"
$q3.Cells.Item(3, 8).Value = 10

$q3.Cells.Item(4, 1).Value = 2
Set-TextValue $q3.Cells.Item(4, 2) "005167"
Set-TextValue $q3.Cells.Item(4, 3) "嘉实润泽量化一年定期开放混合"
Set-TextValue $q3.Cells.Item(4, 4) "0.55"
Set-TextValue $q3.Cells.Item(4, 5) "24.55"
Set-TextValue $q3.Cells.Item(4, 6) "0.67"
Set-TextValue $q3.Cells.Item(4, 7) "0.0037"
$q3.Cells.Item(4, 8).Value = 1

# ---------------------------------------------------------------------------
# Step 4: restyle the header row and column-A cells on "2022-Q3" so they use
# the same bordered/centered format as the "总计" sheet's header (style
# index 2 in the original workbook), matching the format already used for
# this sheet's header/index cells, rather than the old "2022-Q2" format they
# inherited from the Copy() in step 1.
# ---------------------------------------------------------------------------
$total.Cells.Item(1, 2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)
